$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 30300
$ws.Range("C3").Value = 60700
$ws.Range("D3").Value = 90200
$ws.Range("E3").Value = 121000
$ws.Range("B4").Value = 124
$ws.Range("C4").Value = 249
$ws.Range("D4").Value = 369
$ws.Range("E4").Value = 495
$ws.Range("B5").Value = 32.82065
$ws.Range("C5").Value = 32.77107
$ws.Range("D5").Value = 33.04954
$ws.Range("E5").Value = 32.85195
$ws.Range("B11").Value = 9061
$ws.Range("C11").Value = 16500
$ws.Range("D11").Value = 24200
$ws.Range("E11").Value = 29700
$ws.Range("B12").Value = 1188
$ws.Range("C12").Value = 2165
$ws.Range("D12").Value = 3170
$ws.Range("E12").Value = 3890
$ws.Range("B13").Value = 110.08
$ws.Range("C13").Value = 119.57
$ws.Range("D13").Value = 122.19
$ws.Range("E13").Value = 133.41
$ws.Range("B19").Value = 109000
$ws.Range("C19").Value = 178000
$ws.Range("D19").Value = 234000
$ws.Range("E19").Value = 277000
$ws.Range("B20").Value = 444
$ws.Range("C20").Value = 729
$ws.Range("D20").Value = 956
$ws.Range("E20").Value = 1133
$ws.Range("B21").Value = 7.66
$ws.Range("C21").Value = 8.460000000000001
$ws.Range("D21").Value = 8.81
$ws.Range("E21").Value = 9.66
$ws.Range("B27").Value = 9941
$ws.Range("C27").Value = 13300
$ws.Range("D27").Value = 16700
$ws.Range("E27").Value = 17600
$ws.Range("B28").Value = 1303
$ws.Range("C28").Value = 1743
$ws.Range("D28").Value = 2188
$ws.Range("E28").Value = 2304
$ws.Range("B29").Value = 58.04
$ws.Range("C29").Value = 61.17
$ws.Range("D29").Value = 67.8
$ws.Range("E29").Value = 71.55
$ws.Range("B35").Value = 596000
$ws.Range("C35").Value = 728000
$ws.Range("D35").Value = 847000
$ws.Range("E35").Value = 936000
$ws.Range("B36").Value = 2440
$ws.Range("C36").Value = 2983
$ws.Range("D36").Value = 3471
$ws.Range("E36").Value = 3835
$ws.Range("B37").Value = 1.55125
$ws.Range("C37").Value = 2.52192
$ws.Range("D37").Value = 3.30118
$ws.Range("E37").Value = 4.102600000000001
$ws.Range("B43").Value = 15100
$ws.Range("C43").Value = 27700
$ws.Range("D43").Value = 30700
$ws.Range("E43").Value = 32000
$ws.Range("B44").Value = 1974
$ws.Range("C44").Value = 3628
$ws.Range("D44").Value = 4027
$ws.Range("E44").Value = 4194
$ws.Range("B45").Value = 66.19
$ws.Range("C45").Value = 70.42
$ws.Range("D45").Value = 94.93000000000001
$ws.Range("E45").Value = 123.22
$ws.Range("B51").Value = 101000
$ws.Range("C51").Value = 169000
$ws.Range("D51").Value = 216000
$ws.Range("E51").Value = 251000
$ws.Range("B52").Value = 413
$ws.Range("C52").Value = 694
$ws.Range("D52").Value = 883
$ws.Range("E52").Value = 1028
$ws.Range("B53").Value = 7.77
$ws.Range("C53").Value = 8.01
$ws.Range("D53").Value = 8.460000000000001
$ws.Range("E53").Value = 8.640000000000001
$ws.Range("B59").Value = 8062
$ws.Range("C59").Value = 10600
$ws.Range("D59").Value = 12300
$ws.Range("E59").Value = 13400
$ws.Range("B60").Value = 1057
$ws.Range("C60").Value = 1391
$ws.Range("D60").Value = 1617
$ws.Range("E60").Value = 1760
$ws.Range("B61").Value = 64.83
$ws.Range("C61").Value = 66.45
$ws.Range("D61").Value = 71.13
$ws.Range("E61").Value = 74.75
$ws.Range("B67").Value = 30700
$ws.Range("C67").Value = 61200
$ws.Range("D67").Value = 70100
$ws.Range("E67").Value = 118000
$ws.Range("B68").Value = 126
$ws.Range("C68").Value = 251
$ws.Range("D68").Value = 287
$ws.Range("E68").Value = 482
$ws.Range("B69").Value = 32.45162
$ws.Range("C69").Value = 32.49486
$ws.Range("D69").Value = 42.25508
$ws.Range("E69").Value = 33.73512
$ws.Range("B75").Value = 8982
$ws.Range("C75").Value = 17100
$ws.Range("D75").Value = 23800
$ws.Range("E75").Value = 29000
$ws.Range("B76").Value = 1177
$ws.Range("C76").Value = 2237
$ws.Range("D76").Value = 3121
$ws.Range("E76").Value = 3808
$ws.Range("B77").Value = 110.49
$ws.Range("C77").Value = 115.22
$ws.Range("D77").Value = 124.23
$ws.Range("E77").Value = 136.53
$ws.Range("B83").Value = 112000
$ws.Range("C83").Value = 184000
$ws.Range("D83").Value = 241000
$ws.Range("E83").Value = 288000
$ws.Range("B84").Value = 460
$ws.Range("C84").Value = 754
$ws.Range("D84").Value = 987
$ws.Range("E84").Value = 1180
$ws.Range("B85").Value = 7.399850000000001
$ws.Range("C85").Value = 8.242299999999998
$ws.Range("D85").Value = 8.81
$ws.Range("E85").Value = 9.31
$ws.Range("B91").Value = 9941
$ws.Range("C91").Value = 14600
$ws.Range("D91").Value = 16900
$ws.Range("E91").Value = 18100
$ws.Range("B92").Value = 1303
$ws.Range("C92").Value = 1917
$ws.Range("D92").Value = 2212
$ws.Range("E92").Value = 2376
$ws.Range("B93").Value = 58.18084
$ws.Range("C93").Value = 61.38
$ws.Range("D93").Value = 67.64
$ws.Range("E93").Value = 73.06999999999999
$ws.Range("B99").Value = 607000
$ws.Range("C99").Value = 809000
$ws.Range("D99").Value = 855000
$ws.Range("E99").Value = 851000
$ws.Range("B100").Value = 2486
$ws.Range("C100").Value = 3314
$ws.Range("D100").Value = 3501
$ws.Range("E100").Value = 3486
$ws.Range("B101").Value = 1.52803
$ws.Range("C101").Value = 2.34999
$ws.Range("D101").Value = 3.34637
$ws.Range("E101").Value = 4.50437
$ws.Range("B107").Value = 17100
$ws.Range("C107").Value = 25600
$ws.Range("D107").Value = 29300
$ws.Range("E107").Value = 31800
$ws.Range("B108").Value = 2237
$ws.Range("C108").Value = 3355
$ws.Range("D108").Value = 3835
$ws.Range("E108").Value = 4162
$ws.Range("B109").Value = 58.23
$ws.Range("C109").Value = 77.51000000000001
$ws.Range("D109").Value = 98.36
$ws.Range("E109").Value = 124.42
$ws.Range("B115").Value = 101000
$ws.Range("C115").Value = 171000
$ws.Range("D115").Value = 217000
$ws.Range("E115").Value = 239000
$ws.Range("B116").Value = 416
$ws.Range("C116").Value = 701
$ws.Range("D116").Value = 889
$ws.Range("E116").Value = 980
$ws.Range("B117").Value = 7.93
$ws.Range("C117").Value = 7.974810000000001
$ws.Range("D117").Value = 8.470000000000001
$ws.Range("E117").Value = 9.300000000000001
$ws.Range("B123").Value = 7937
$ws.Range("C123").Value = 9061
$ws.Range("D123").Value = 12100
$ws.Range("E123").Value = 12900
$ws.Range("B124").Value = 1040
$ws.Range("C124").Value = 1188
$ws.Range("D124").Value = 1592
$ws.Range("E124").Value = 1688
$ws.Range("B125").Value = 65.61114000000001
$ws.Range("C125").Value = 102.56
$ws.Range("D125").Value = 79.44
$ws.Range("E125").Value = 78.59999999999999
